$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.63
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.4
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("AA5").Value = 23
$ws.Range("AC5").Value = 7
$ws.Range("AH5").Value = 7.5
$ws.Range("AU5").Value = 8.5
$ws.Range("BA5").Value = 81
$ws.Range("BB5").Value = 251

$ws.Range("G24").Value = 3.1
$ws.Range("I24").Value = 2.35
$ws.Range("J24").Value = 3.75
$ws.Range("L24").Value = 3.1
$ws.Range("M24").Value = 1.07
$ws.Range("N24").Value = 9
$ws.Range("O24").Value = 1.36
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = 2.2
$ws.Range("R24").Value = 1.65
$ws.Range("W24").Value = 8.5
$ws.Range("X24").Value = 15
$ws.Range("Y24").Value = 12
$ws.Range("Z24").Value = 34
$ws.Range("AA24").Value = 29
$ws.Range("AI24").Value = 10
$ws.Range("AK24").Value = 21
$ws.Range("AL24").Value = 21
$ws.Range("AN24").Value = 5
$ws.Range("AO24").Value = 19
$ws.Range("AW24").Value = 4.33
$ws.Range("AX24").Value = 13
$ws.Range("AZ24").Value = 41
$ws.Range("BA24").Value = 67
